# Daily attendance processing - 2025-10-19 21:16:58
# Normalizes the "Recorded By" (column G) values so that the literal
# "System" token is moved from the front of the comma-separated list
# to the end, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Text

    if ($null -ne $val -and $val.StartsWith("System, ")) {
        $parts = $val -split ", "
        if ($parts[0] -eq "System") {
            $rest = $parts[1..($parts.Length - 1)]
            $newVal = ($rest -join ", ") + ", System"
            $cell.Value = $newVal
        }
    }
}
